$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format columns I and L as Text so date-like strings ("dd/mm/yyyy")
# are stored as literal text (matching the source data) instead of being
# auto-converted into Excel date serials.
$ws.Range("I442").NumberFormat = "@"
$ws.Range("L442").NumberFormat = "@"

# Row 433
$ws.Cells.Item(433,1).Value = "Gessé Rourera"
$ws.Cells.Item(433,2).Value = 22
$ws.Cells.Item(433,3).Value = 1849
$ws.Cells.Item(433,16).Value = "3 Baptismes 1831-1851"

# Row 434
$ws.Cells.Item(434,1).Value = "Gras Ros"
$ws.Cells.Item(434,2).Value = 22
$ws.Cells.Item(434,3).Value = 1849
$ws.Cells.Item(434,16).Value = "3 Baptismes 1831-1851"

# Row 435
$ws.Cells.Item(435,1).Value = "Castells Bendicho"
$ws.Cells.Item(435,2).Value = 22
$ws.Cells.Item(435,3).Value = 1849
$ws.Cells.Item(435,16).Value = "3 Baptismes 1831-1851"

# Row 436
$ws.Cells.Item(436,1).Value = "Pedra Gessé"
$ws.Cells.Item(436,2).Value = 22
$ws.Cells.Item(436,3).Value = 1850
$ws.Cells.Item(436,16).Value = "3 Baptismes 1831-1851"

# Row 437
$ws.Cells.Item(437,1).Value = "Gessé Queralt"
$ws.Cells.Item(437,2).Value = 23
$ws.Cells.Item(437,3).Value = 1850
$ws.Cells.Item(437,16).Value = "3 Baptismes 1831-1851"

# Row 438
$ws.Cells.Item(438,1).Value = "Rodrigo Gessé"
$ws.Cells.Item(438,2).Value = 23
$ws.Cells.Item(438,3).Value = 1850
$ws.Cells.Item(438,16).Value = "3 Baptismes 1831-1851"

# Row 439
$ws.Cells.Item(439,1).Value = "Gessé Taribó"
$ws.Cells.Item(439,2).Value = 23
$ws.Cells.Item(439,3).Value = 1850
$ws.Cells.Item(439,16).Value = "3 Baptismes 1831-1851"

# Row 440
$ws.Cells.Item(440,1).Value = "Domingo Pedra"
$ws.Cells.Item(440,2).Value = 23
$ws.Cells.Item(440,3).Value = 1850
$ws.Cells.Item(440,16).Value = "3 Baptismes 1831-1851"

# Row 441
$ws.Cells.Item(441,1).Value = "Gessé Mauri"
$ws.Cells.Item(441,2).Value = 23
$ws.Cells.Item(441,3).Value = 1851
$ws.Cells.Item(441,16).Value = "3 Baptismes 1831-1851"

# Row 442
$ws.Cells.Item(442,1).Value = "Gessé Gessé Joseph Liberato Ramon"
$ws.Cells.Item(442,2).Value = 24
$ws.Cells.Item(442,3).Value = 1851
$ws.Cells.Item(442,4).Value = "Jaume Gessé Marquet"
$ws.Cells.Item(442,5).Value = "Vicenta Gessé Marquet"
$ws.Cells.Item(442,6).Value = "Joseph Gessé i Rita Marquet"
$ws.Cells.Item(442,7).Value = "Thomas Gessé i Rita Marquet"
$ws.Cells.Item(442,9).Value = "06/05/1851"
$ws.Cells.Item(442,10).Value = "Joseph Gessé"
$ws.Cells.Item(442,11).Value = "Rosa Cortés"
$ws.Cells.Item(442,12).Value = "06/05/1851"
$ws.Cells.Item(442,16).Value = "3 Baptismes 1831-1851"

# Row 443
$ws.Cells.Item(443,1).Value = "Puig Gessé"
$ws.Cells.Item(443,2).Value = 24
$ws.Cells.Item(443,3).Value = 1851
$ws.Cells.Item(443,16).Value = "3 Baptismes 1831-1851"

# Row 444
$ws.Cells.Item(444,1).Value = "Pedra Gessé"
$ws.Cells.Item(444,2).Value = 24
$ws.Cells.Item(444,3).Value = 1851
$ws.Cells.Item(444,16).Value = "3 Baptismes 1831-1851"

# Row 445
$ws.Cells.Item(445,1).Value = "Badia Estada"
$ws.Cells.Item(445,2).Value = 25
$ws.Cells.Item(445,3).Value = 1851
$ws.Cells.Item(445,16).Value = "3 Baptismes 1831-1851"

# Row 446
$ws.Cells.Item(446,1).Value = "Gessé Badia"
$ws.Cells.Item(446,2).Value = 25
$ws.Cells.Item(446,3).Value = 1851
$ws.Cells.Item(446,16).Value = "3 Baptismes 1831-1851"

# Row 447
$ws.Cells.Item(447,1).Value = "Nadal Gessé"
$ws.Cells.Item(447,2).Value = 25
$ws.Cells.Item(447,3).Value = 1851
$ws.Cells.Item(447,16).Value = "3 Baptismes 1831-1851"

# Row 448
$ws.Cells.Item(448,1).Value = "Badia Castells"
$ws.Cells.Item(448,2).Value = 26
$ws.Cells.Item(448,3).Value = 1851
$ws.Cells.Item(448,16).Value = "3 Baptismes 1831-1851"

# Row 449
$ws.Cells.Item(449,1).Value = "Badia Cortés"
$ws.Cells.Item(449,2).Value = 26
$ws.Cells.Item(449,3).Value = 1852
$ws.Cells.Item(449,16).Value = "3 Baptismes 1831-1851"

# Row 450
$ws.Cells.Item(450,1).Value = "Seufenís?"
$ws.Cells.Item(450,2).Value = 1
$ws.Cells.Item(450,3).Value = 1737
$ws.Cells.Item(450,16).Value = "2 Quinque Libri 1750-1831"

# Row 451
$ws.Cells.Item(451,1).Value = "Jubillà"
$ws.Cells.Item(451,2).Value = 1
$ws.Cells.Item(451,3).Value = 1737

# Row 452
$ws.Cells.Item(452,1).Value = "Terés"
$ws.Cells.Item(452,2).Value = 1
$ws.Cells.Item(452,3).Value = 1737

# Row 453
$ws.Cells.Item(453,1).Value = "Gessé"
$ws.Cells.Item(453,2).Value = 2
$ws.Cells.Item(453,3).Value = 1738

# Row 454
$ws.Cells.Item(454,1).Value = "Castells"
$ws.Cells.Item(454,2).Value = 2
$ws.Cells.Item(454,3).Value = 1738

# Row 455
$ws.Cells.Item(455,1).Value = "Gessé"
$ws.Cells.Item(455,2).Value = 2
$ws.Cells.Item(455,3).Value = 1739

# Row 456
$ws.Cells.Item(456,1).Value = "Jubillà"
$ws.Cells.Item(456,2).Value = 3
$ws.Cells.Item(456,3).Value = 1739

# Row 457
$ws.Cells.Item(457,1).Value = "Gessé"
$ws.Cells.Item(457,2).Value = 3
$ws.Cells.Item(457,3).Value = 1740

# Row 458
$ws.Cells.Item(458,1).Value = "Castells"
$ws.Cells.Item(458,2).Value = 3
$ws.Cells.Item(458,3).Value = 1740

# Row 459
$ws.Cells.Item(459,1).Value = "Badia"
$ws.Cells.Item(459,2).Value = 4
$ws.Cells.Item(459,3).Value = 1740

# Row 460
$ws.Cells.Item(460,1).Value = "Terés"
$ws.Cells.Item(460,2).Value = 4
$ws.Cells.Item(460,3).Value = 1740

# Row 461
$ws.Cells.Item(461,1).Value = "Castells"
$ws.Cells.Item(461,2).Value = 4
$ws.Cells.Item(461,3).Value = 1741

# Row 462
$ws.Cells.Item(462,1).Value = "Badia"
$ws.Cells.Item(462,2).Value = 5
$ws.Cells.Item(462,3).Value = 1742

# Row 463
$ws.Cells.Item(463,1).Value = "Jubillà"
$ws.Cells.Item(463,2).Value = 5
$ws.Cells.Item(463,3).Value = 1742

# Row 464
$ws.Cells.Item(464,1).Value = "Gessé Amat"
$ws.Cells.Item(464,2).Value = 5
$ws.Cells.Item(464,3).Value = 1743

# Row 465
$ws.Cells.Item(465,1).Value = "Ferrer Molins"
$ws.Cells.Item(465,2).Value = 6
$ws.Cells.Item(465,3).Value = 1743

# Row 466
$ws.Cells.Item(466,1).Value = "Badia Oliva"
$ws.Cells.Item(466,2).Value = 6
$ws.Cells.Item(466,3).Value = 1744

# Row 467
$ws.Cells.Item(467,1).Value = "Castells Mas"
$ws.Cells.Item(467,2).Value = 6
$ws.Cells.Item(467,3).Value = 1744

# Row 468
$ws.Cells.Item(468,1).Value = "Jubillà"
$ws.Cells.Item(468,2).Value = 7
$ws.Cells.Item(468,3).Value = 1745

# Row 469
$ws.Cells.Item(469,1).Value = "Porta"
$ws.Cells.Item(469,2).Value = 7
$ws.Cells.Item(469,3).Value = 1745

# Row 470
$ws.Cells.Item(470,1).Value = "Castells Mas"
$ws.Cells.Item(470,2).Value = 8
$ws.Cells.Item(470,3).Value = 1746

# Row 471
$ws.Cells.Item(471,1).Value = "Farré"
$ws.Cells.Item(471,2).Value = 8
$ws.Cells.Item(471,3).Value = 1746

# Row 472
$ws.Cells.Item(472,1).Value = "Farré"
$ws.Cells.Item(472,2).Value = 8
$ws.Cells.Item(472,3).Value = 1746

# Row 473
$ws.Cells.Item(473,1).Value = "Castells Mas"
$ws.Cells.Item(473,2).Value = 9
$ws.Cells.Item(473,3).Value = 1747

# Row 474
$ws.Cells.Item(474,1).Value = "Porta"
$ws.Cells.Item(474,2).Value = 9
$ws.Cells.Item(474,3).Value = 1748

# Row 475
$ws.Cells.Item(475,1).Value = "Molins Gessé"
$ws.Cells.Item(475,2).Value = 10
$ws.Cells.Item(475,3).Value = 1748

# Row 476
$ws.Cells.Item(476,1).Value = "Castells Hospital"
$ws.Cells.Item(476,2).Value = 10
$ws.Cells.Item(476,3).Value = 1748

# Row 477
$ws.Cells.Item(477,1).Value = "Terés Flores"
$ws.Cells.Item(477,2).Value = 10
$ws.Cells.Item(477,3).Value = 1748

# Update the active selection to match the new last row (A478), as Excel
# would leave it positioned after the last data-entry row.
$ws.Range("A478").Select()

